$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 29333.334
$ws.Range("I69").Value = 29142.857
$ws.Range("K69").Value = 87428.571
$ws.Range("M69").Value = -86554.571
$ws.Range("H72").Value = 29333.334
$ws.Range("I72").Value = 29142.857
$ws.Range("K72").Value = 262285.713
$ws.Range("M72").Value = -257917.713
$ws.Range("H98").Value = 1348.25
$ws.Range("I98").Value = 1348.25
$ws.Range("K98").Value = 1348.25
$ws.Range("M98").Value = 149.75
$ws.Range("H100").Value = 2668
$ws.Range("I100").Value = 2309.6
$ws.Range("K100").Value = 2309.6
$ws.Range("M100").Value = -1768.6
$ws.Range("H107").Value = 1352.3077
$ws.Range("I107").Value = 275
$ws.Range("K107").Value = 275
$ws.Range("M107").Value = 1645
$ws.Range("H122").Value = 1348.25
$ws.Range("I122").Value = 1348.25
$ws.Range("K122").Value = 4044.75
$ws.Range("M122").Value = -1594.75
$ws.Range("H135").Value = 948.8889
$ws.Range("I135").Value = 948.8889
$ws.Range("K135").Value = 8540.000100000001
$ws.Range("M135").Value = -6005.000100000001
$ws.Range("H137").Value = 559558.75
$ws.Range("I137").Value = 1622.1428
$ws.Range("J137").Value = 1210484.9
$ws.Range("K137").Value = 4866.428400000001
$ws.Range("L137").Value = 3631454.7
$ws.Range("M137").Value = -2316.428400000001
$ws.Range("N137").Value = -3636554.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2975.2307
$ws.Range("I45").Value = 2852.5454
$ws.Range("K45").Value = 2852.5454
$ws.Range("M45").Value = -2475.5454
$ws.Range("H61").Value = 2116.15
$ws.Range("I61").Value = 1871.9412
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 1871.9412
$ws.Range("L61").Value = 3500
$ws.Range("M61").Value = -1659.9412
$ws.Range("N61").Value = -3924
$ws.Range("H92").Value = 14999
$ws.Range("J92").Value = 14999
$ws.Range("L92").Value = 14999
$ws.Range("N92").Value = -19991
$ws.Range("H96").Value = 50000
$ws.Range("J96").Value = 50000
$ws.Range("L96").Value = 50000
$ws.Range("N96").Value = -55492
$ws.Range("H102").Value = 53850.91
$ws.Range("I102").Value = 54287.316
$ws.Range("K102").Value = 54287.316
$ws.Range("M102").Value = -52665.316
$ws.Range("H110").Value = 1310.8125
$ws.Range("I110").Value = 1180.7142
$ws.Range("J110").Value = 2221.5
$ws.Range("K110").Value = 1180.7142
$ws.Range("L110").Value = 2221.5
$ws.Range("M110").Value = 864.2858000000001
$ws.Range("N110").Value = -6311.5
$ws.Range("H122").Value = 3945.0435
$ws.Range("I122").Value = 4399.4116
$ws.Range("J122").Value = 2657.6667
$ws.Range("K122").Value = 13198.2348
$ws.Range("L122").Value = 7973.000100000001
$ws.Range("M122").Value = -10748.2348
$ws.Range("N122").Value = -12873.0001
$ws.Range("H125").Value = 56347.332
$ws.Range("J125").Value = 56347.332
$ws.Range("L125").Value = 56347.332
$ws.Range("N125").Value = -66187.33199999999
$ws.Range("H132").Value = 2508.76
$ws.Range("I132").Value = 1819.4667
$ws.Range("K132").Value = 5458.4001
$ws.Range("M132").Value = -2928.4001
$ws.Range("H136").Value = 2116.15
$ws.Range("I136").Value = 1871.9412
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 5615.8236
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -3065.8236
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 97910.13
$ws.Range("I20").Value = 124209.11
$ws.Range("K20").Value = 124209.11
$ws.Range("M20").Value = -123962.11
$ws.Range("H80").Value = 62942.125
$ws.Range("I80").Value = 250268.75
$ws.Range("K80").Value = 250268.75
$ws.Range("M80").Value = -249270.75
$ws.Range("H83").Value = 62942.125
$ws.Range("I83").Value = 250268.75
$ws.Range("K83").Value = 1251343.75
$ws.Range("M83").Value = -1246351.75
$ws.Range("H86").Value = 5901.6
$ws.Range("I86").Value = 5747.5
$ws.Range("J86").Value = 6004.3335
$ws.Range("K86").Value = 5747.5
$ws.Range("L86").Value = 6004.3335
$ws.Range("M86").Value = -4624.5
$ws.Range("N86").Value = -8250.333500000001
$ws.Range("H89").Value = 5901.6
$ws.Range("I89").Value = 5747.5
$ws.Range("J89").Value = 6004.3335
$ws.Range("K89").Value = 28737.5
$ws.Range("L89").Value = 30021.6675
$ws.Range("M89").Value = -23121.5
$ws.Range("N89").Value = -41253.6675
$ws.Range("H107").Value = 4660.143
$ws.Range("I107").Value = 3957.75
$ws.Range("J107").Value = 5596.6665
$ws.Range("K107").Value = 3957.75
$ws.Range("L107").Value = 5596.6665
$ws.Range("M107").Value = -2037.75
$ws.Range("N107").Value = -9436.666499999999
$ws.Range("H134").Value = 2315.261
$ws.Range("I134").Value = 1678.3125
$ws.Range("J134").Value = 3771.1428
$ws.Range("K134").Value = 5034.9375
$ws.Range("L134").Value = 11313.4284
$ws.Range("M134").Value = -2499.9375
$ws.Range("N134").Value = -16383.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 1601
$ws.Range("J92").Value = 1601
$ws.Range("L92").Value = 1601
$ws.Range("N92").Value = -6593
$ws.Range("H94").Value = 933.6667
$ws.Range("I94").Value = 906.4545000000001
$ws.Range("J94").Value = 952.375
$ws.Range("K94").Value = 906.4545000000001
$ws.Range("L94").Value = 952.375
$ws.Range("M94").Value = -455.4545000000001
$ws.Range("N94").Value = -1854.375
$ws.Range("H96").Value = 35812
$ws.Range("J96").Value = 35812
$ws.Range("L96").Value = 35812
$ws.Range("N96").Value = -41304
$ws.Range("H103").Value = 1524
$ws.Range("I103").Value = 1524
$ws.Range("K103").Value = 1524
$ws.Range("M103").Value = -352
$ws.Range("H132").Value = 2709.6667
$ws.Range("J132").Value = 3187
$ws.Range("L132").Value = 9561
$ws.Range("N132").Value = -14621
$ws.Range("H134").Value = 2192.6956
$ws.Range("I134").Value = 1864.4286
$ws.Range("K134").Value = 5593.2858
$ws.Range("M134").Value = -3058.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 625381.0600000001
$ws.Range("I97").Value = 833643.8
$ws.Range("J97").Value = 592.75
$ws.Range("K97").Value = 833643.8
$ws.Range("L97").Value = 592.75
$ws.Range("M97").Value = -833147.8
$ws.Range("N97").Value = -1584.75
$ws.Range("H113").Value = 8612964
$ws.Range("I113").Value = 1111011
$ws.Range("K113").Value = 1111011
$ws.Range("M113").Value = -1108841
$ws.Range("H132").Value = 8486.272000000001
$ws.Range("I132").Value = 8514.143
$ws.Range("K132").Value = 25542.429
$ws.Range("M132").Value = -23012.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3064.5
$ws.Range("J46").Value = 3064.5
$ws.Range("L46").Value = 3064.5
$ws.Range("N46").Value = -3440.5
$ws.Range("H132").Value = 10627.723
$ws.Range("I132").Value = 17320.889
$ws.Range("K132").Value = 51962.667
$ws.Range("M132").Value = -49432.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 50000
$ws.Range("I40").Value = 50000
$ws.Range("K40").Value = 50000
$ws.Range("M40").Value = -49851
$ws.Range("H41").Value = 210187.5
$ws.Range("J41").Value = 210187.5
$ws.Range("L41").Value = 210187.5
$ws.Range("N41").Value = -210967.5
$ws.Range("H122").Value = 4421.1113
$ws.Range("I122").Value = 3338
$ws.Range("K122").Value = 10014
$ws.Range("M122").Value = -7564
